$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price/volume snapshot.
# Values are prefixed with a leading apostrophe (PowerShell single-
# quoted string so it is written literally) which Excel treats as an
# explicit "store as text" marker. Without it, numeric-looking strings
# such as "1.000" or "0.000009057" get coerced to Double and lose their
# original formatting, which the source data (and the diff) relies on.

$ws.Range('D2').Value = '''29.027.88'
$ws.Range('E2').Value = '''  -0.41%  '
$ws.Range('D3').Value = '''1.829.32'
$ws.Range('E3').Value = '''  -0.07%  '
$ws.Range('D4').Value = '''0.9990'
$ws.Range('E4').Value = '''  -0.02%  '
$ws.Range('E5').Value = '''  -0.26%  '
$ws.Range('D6').Value = '''0.6236'
$ws.Range('E6').Value = '''  -5.20%  '
$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '''  +0.01%  '
$ws.Range('B8').Value = '''Dogecoin'
$ws.Range('C8').Value = '''https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').Value = '''0.07527'
$ws.Range('E8').Value = '''  +1.83%  '
$ws.Range('B9').Value = '''OKB'
$ws.Range('C9').Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = '''44.62'
$ws.Range('E9').Value = '''  +6.74%  '
$ws.Range('D10').Value = '''0.2908'
$ws.Range('E10').Value = '''  -0.41%  '
$ws.Range('D11').Value = '''22.76'
$ws.Range('E11').Value = '''  -0.29%  '
$ws.Range('D12').Value = '''0.07642'
$ws.Range('E12').Value = '''  -1.40%  '
$ws.Range('D13').Value = '''1.829.56'
$ws.Range('E13').Value = '''  -0.50%  '
$ws.Range('D14').Value = '''4.957'
$ws.Range('E14').Value = '''  -0.60%  '
$ws.Range('D15').Value = '''0.6638'
$ws.Range('E15').Value = '''  -0.10%  '
$ws.Range('D16').Value = '''82.27'
$ws.Range('E16').Value = '''  -0.54%  '
$ws.Range('D17').Value = '''0.000009057'
$ws.Range('E17').Value = '''  +7.55%  '
$ws.Range('D18').Value = '''6.004'
$ws.Range('E18').Value = '''  -1.94%  '
$ws.Range('D19').Value = '''28.923.91'
$ws.Range('E19').Value = '''  -0.81%  '
$ws.Range('D20').Value = '''224.76'
$ws.Range('E20').Value = '''  -0.86%  '
$ws.Range('E21').Value = '''  -0.83%  '
$ws.Range('D22').Value = '''1.001'
$ws.Range('E22').Value = '''  +0.02%  '
$ws.Range('D23').Value = '''7.188'
$ws.Range('E23').Value = '''  +1.06%  '
$ws.Range('E24').Value = '''  +0.07%  '
$ws.Range('D25').Value = '''159.30'
$ws.Range('E25').Value = '''  +0.52%  '
$ws.Range('D26').Value = '''8.379'
$ws.Range('E26').Value = '''  -2.42%  '
$ws.Range('D27').Value = '''0.1355'
$ws.Range('E27').Value = '''  -2.55%  '
$ws.Range('D28').Value = '''17.82'
$ws.Range('E28').Value = '''  -0.38%  '
$ws.Range('D29').Value = '''1.493'
$ws.Range('E29').Value = '''  -1.62%  '
$ws.Range('D30').Value = '''4.033'
$ws.Range('E30').Value = '''  -0.17%  '
$ws.Range('D31').Value = '''4.050'
$ws.Range('E31').Value = '''  -1.40%  '
$ws.Range('D32').Value = '''1.200'
$ws.Range('E32').Value = '''  +0.84%  '
$ws.Range('D33').Value = '''0.05199'
$ws.Range('E33').Value = '''  -0.98%  '
$ws.Range('D34').Value = '''1.837'
$ws.Range('E34').Value = '''  -1.38%  '
$ws.Range('D35').Value = '''1.152'
$ws.Range('E35').Value = '''  +1.03%  '
$ws.Range('D36').Value = '''0.7316'
$ws.Range('E36').Value = '''  -0.95%  '
$ws.Range('D37').Value = '''2.609'
$ws.Range('E37').Value = '''  -1.64%  '
$ws.Range('D38').Value = '''1.285.08'
$ws.Range('E38').Value = '''  -1.12%  '
$ws.Range('D39').Value = '''2.758'
$ws.Range('E39').Value = '''  +0.98%  '
$ws.Range('D40').Value = '''0.01779'
$ws.Range('E40').Value = '''  -0.66%  '
$ws.Range('D41').Value = '''6.383'
$ws.Range('E41').Value = '''  +7.23%  '
$ws.Range('D42').Value = '''0.8948'
$ws.Range('E42').Value = '''  -2.65%  '
$ws.Range('D44').Value = '''101.25'
$ws.Range('E44').Value = '''  -0.82%  '
$ws.Range('D45').Value = '''1.980.24'
$ws.Range('E45').Value = '''  +0.19%  '
$ws.Range('D46').Value = '''0.5115'
$ws.Range('E46').Value = '''  -0.46%  '
$ws.Range('D47').Value = '''63.47'
$ws.Range('E47').Value = '''  +0.42%  '
$ws.Range('E48').Value = '''  -0.49%  '
$ws.Range('D49').Value = '''0.3971'
$ws.Range('E49').Value = '''  -0.70%  '
$ws.Range('D50').Value = '''8.854'
$ws.Range('E50').Value = '''  +1.26%  '
$ws.Range('D51').Value = '''1.652'
$ws.Range('E51').Value = '''  -5.57%  '
